# LeetCode Practice Tracker update
# - Logs two newly solved Sliding Window problems (rows 62-63):
#     "Subarrays with K Different Integers" and "Binary Subarrays With Sum"
# - Jots down three more problem names to tackle later (rows 64-66)
# - Pushes the previously queued "To Do" rows further down the sheet
#   (rows 64-68 -> rows 70-74) to make room for the new entries above

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: move the existing queued rows further down ---
# Topic (col B) for the two rows that actually change identity
$ws.Range("B66:B68").Cut($ws.Range("B72:B74"))
$excel.CutCopyMode = $false

# Problem name / difficulty / status (cols C:E) for all five rows
$ws.Range("C64:E68").Cut($ws.Range("C70:E74"))
$excel.CutCopyMode = $false

# Clean up the blank cell stubs the cut left behind in the vacated rows
$ws.Range("B66:E68").Clear()
$ws.Range("C64:E65").Clear()

# Rows 70/71 keep the same "Sliding Window" topic the rows had before the move
$ws.Range("B70").Value = "Sliding Window"
$ws.Range("B71").Value = "Sliding Window"

# --- Row 62: Subarrays with K Different Integers (solved) ---
$ws.Range("B62").Value = "Sliding Window"
$ws.Range("C62").Value = "Subarrays with K Different Integers"
$ws.Range("D62").Value = "Medium"
$ws.Range("E62").Value = "Done"
$ws.Range("F61").Copy()
$ws.Range("F62:F63").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F62").Value = 45901
$ws.Range("G62").Value = "O(n)"
$ws.Range("H62").Value = "O(n)"
$ws.Range("I62").Value = "Sliding Window + Hmap"

# --- Row 63: Binary Subarrays With Sum (solved) ---
$ws.Range("B63").Value = "Sliding Window"
$ws.Range("C63").Value = "Binary Subarrays With Sum"
$ws.Range("D63").Value = "Medium"
$ws.Range("E63").Value = "Done"
$ws.Range("F63").Value = 45901
$ws.Range("G63").Value = "O(n)"
$ws.Range("H63").Value = "O(1)"
$ws.Range("I63").Value = "Sliding Window"

# --- Rows 64-66: problem names noted for later ---
$ws.Range("C64").Value = "Minimum Window Substring"
$ws.Range("C65").Value = "Sliding Window Maximum"
$ws.Range("C66").Value = "Longest Substring with At Least K Repeating Characters"

# --- Restore view state (scroll position / selection) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("H64").Select()

Write-Host "Applied LeetCode tracker update"
